$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "[대학원] 내용은 둘째고, 우리랑 걔네들이랑 쓰는 언어가 다르잖아"
$ws.Range("E9").Value = "https://pdsi.pabii.com/different-language-different-context/#utm_source=rss&utm_medium=rss&utm_campaign=different-language-different-context"

$ws.Range("D27").Value = "Luda Gen 1, 더 재미있고 자연스러운 대화로 돌아온 루다 2편 - 생성 모델을 챗봇으로 빚어내기"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/luda-gen-2/"

$ws.Range("D51").Value = "[python] 타입 힌트와 관련하여 중요한 모듈인 typing을 알아보자"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%ED%83%80%EC%9E%85-%ED%9E%8C%ED%8A%B8%EC%99%80-%EA%B4%80%EB%A0%A8%ED%95%98%EC%97%AC-%EC%A4%91%EC%9A%94%ED%95%9C-%EB%AA%A8%EB%93%88%EC%9D%B8-typing%EC%9D%84-%EC%95%8C%EC%95%84%EB%B3%B4%EC%9E%90"
